# Fix POC pCO2 extraction mismatch
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (rows 2..13, columns A..E) per the corrected ascertainment counts
$data = @(
    @($false, $false, $true,  $false, 8694),
    @($false, $false, $false, $true,  8409),
    @($false, $false, $true,  $true,  5059),
    @($false, $true,  $false, $true,  2142),
    @($false, $true,  $true,  $true,  1670),
    @($true,  $true,  $true,  $true,  701),
    @($true,  $false, $true,  $false, 438),
    @($true,  $false, $true,  $true,  326),
    @($true,  $false, $false, $false, 316),
    @($true,  $true,  $false, $true,  129),
    @($true,  $false, $false, $true,  73),
    @($false, $true,  $true,  $false, 18)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
    $ws.Range("D$row").Value = $data[$i][3]
    $ws.Range("E$row").Value = $data[$i][4]
}
